$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest refreshed values. Price cells whose new text would
# otherwise be auto-parsed by Excel as a number (e.g. "574.40") are
# written with a leading apostrophe so they stay plain text, matching
# the original inline-string cell contents.
$ws.Range('D2').Value = '63.900.71'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.751.30'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''574.40'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').Value = '''157.36'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('E9').Value = '  -3.24%  '
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('E11').Value = '  -16.57%  '
$ws.Range('E12').Value = '  -3.35%  '
$ws.Range('D13').Value = '3.238.66'
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('D14').Value = '''26.48'
$ws.Range('E14').Value = '  -3.29%  '
$ws.Range('D15').Value = '63.657.95'
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('D17').Value = '2.754.60'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('E19').Value = '  -2.98%  '
$ws.Range('D20').Value = '''354.54'
$ws.Range('E20').Value = '  -2.37%  '
$ws.Range('E21').Value = '  -4.27%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '''0.533'
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').Value = '''65.13'
$ws.Range('E24').Value = '  -2.97%  '
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').Value = '''8.45'
$ws.Range('E27').Value = '  -2.16%  '
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('E29').Value = '  -4.88%  '
$ws.Range('E30').Value = '  -1.73%  '
$ws.Range('E31').Value = '  -1.19%  '
$ws.Range('D32').Value = '''168.50'
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('E33').Value = '  -2.71%  '
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  -0.44%  '
$ws.Range('E37').Value = '  -1.77%  '
$ws.Range('D38').Value = '''0.984'
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('D39').Value = '''6.16'
$ws.Range('E39').Value = '  +5.74%  '
$ws.Range('D40').Value = '''331.75'
$ws.Range('E40').Value = '  -2.49%  '
$ws.Range('E41').Value = '  -4.26%  '
$ws.Range('D42').Value = '''38.96'
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').Value = '''21.44'
$ws.Range('E43').Value = '  -2.56%  '
$ws.Range('D44').Value = '''0.0587'
$ws.Range('E44').Value = '  -2.35%  '
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('D47').Value = '''134.95'
$ws.Range('E47').Value = '  -2.26%  '
$ws.Range('E48').Value = '  -4.05%  '
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('E50').Value = '  -0.13%  '
$ws.Range('D51').Value = '''11.05'
$ws.Range('E51').Value = '  +0.23%  '
